# Update category classifications in column F (Test Type) and a handful of
# Assessment Name entries in column A on Sheet1, then leave the selection /
# zoom matching where the author ended up after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Assessment Name) corrections ---
$ws.Range("A4").Value   = "Web Development"
$ws.Range("A30").Value  = "Data Science"
$ws.Range("A54").Value  = "Data Science"
$ws.Range("A91").Value  = "Data Science"
$ws.Range("A93").Value  = "Data Science"
$ws.Range("A97").Value  = "Data Science"
$ws.Range("A100").Value = "Data Science"

# --- Column F (Test Type) corrections ---
$ws.Range("F2").Value   = "Knowledge & Skills"
$ws.Range("F3").Value   = "Agile Testing"
$ws.Range("F4").Value   = "Programming"
$ws.Range("F6").Value   = "Simulation"
$ws.Range("F9").Value   = "Programming"
$ws.Range("F12").Value  = "Programming"
$ws.Range("F13").Value  = "Simulation"
$ws.Range("F15").Value  = "Agile Testing"
$ws.Range("F20").Value  = "Knowledge & Skills"
$ws.Range("F23").Value  = "Programming"
$ws.Range("F25").Value  = "Knowledge & Skills"
$ws.Range("F26").Value  = "Knowledge & Skills"
$ws.Range("F46").Value  = "Knowledge & Skills"
$ws.Range("F58").Value  = "Programming & Skills"
$ws.Range("F73").Value  = "Knowledge & Skills"
$ws.Range("F77").Value  = "Knowledge & Skills"
$ws.Range("F83").Value  = "Knowledge & Skills"
$ws.Range("F84").Value  = "Knowledge & Skills"
$ws.Range("F87").Value  = "Knowledge & Skills"
$ws.Range("F88").Value  = "Knowledge & Skills"
$ws.Range("F91").Value  = "Knowledge & Skills"
$ws.Range("F94").Value  = "Knowledge & Skills"
$ws.Range("F100").Value = "Knowledge & Skills"

# --- View state: scroll/zoom/selection to match where the author ended up ---
$excel.ActiveWindow.Zoom = 74
[void]$ws.Range("F101").Select()
